# Weekly update: insert a new data row at row 12 (pushing existing rows
# 12..40 down to 13..41) and populate it with the new week's record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 12..40 down by one row to make room for the new record.
$ws.Rows.Item(12).Insert()

# Populate the newly inserted row 12 with the new data.
$ws.Cells.Item(12, 1).Value = 1
$ws.Cells.Item(12, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(12, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(12, 4).Value = 45036
$ws.Cells.Item(12, 5).Value = 15
$ws.Cells.Item(12, 6).Value = 100112044
$ws.Cells.Item(12, 7).Value = "Perejil"
$ws.Cells.Item(12, 8).Value = "Sin especificar"
$ws.Cells.Item(12, 9).Value = "Segunda"
$ws.Cells.Item(12, 10).Value = 210
$ws.Cells.Item(12, 11).Value = 2300
$ws.Cells.Item(12, 12).Value = 2500
$ws.Cells.Item(12, 13).Value = 2443
$ws.Cells.Item(12, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(12, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(12, 16).Value = 1222
$ws.Cells.Item(12, 17).Value = 2
$ws.Cells.Item(12, 18).Value = "Hortaliza"
